# Auto-applied numeric corrections to the Leve profit-tracking columns
# (currentAveragePrice*, Leve*Price*, Leve*Profit*) across all 8 job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1231.5714
$ws.Range("J17").Value = 1231.5714
$ws.Range("L17").Value = 3694.7142
$ws.Range("N17").Value = -4030.7142

$ws.Range("H19").Value = 801.8095
$ws.Range("I19").Value = 658.7273
$ws.Range("J19").Value = 959.2
$ws.Range("K19").Value = 658.7273
$ws.Range("L19").Value = 959.2
$ws.Range("M19").Value = -483.7273
$ws.Range("N19").Value = -1309.2

$ws.Range("H33").Value = 705.3200000000001
$ws.Range("I33").Value = 696.7727
$ws.Range("J33").Value = 768
$ws.Range("K33").Value = 696.7727
$ws.Range("L33").Value = 768
$ws.Range("M33").Value = -467.7727
$ws.Range("N33").Value = -1226

$ws.Range("H40").Value = 4066.9167
$ws.Range("I40").Value = 6250.5
$ws.Range("J40").Value = 3630.2
$ws.Range("K40").Value = 6250.5
$ws.Range("L40").Value = 3630.2
$ws.Range("M40").Value = -6075.5
$ws.Range("N40").Value = -3980.2

$ws.Range("H132").Value = 4997.1045
$ws.Range("I132").Value = 4747.7437
$ws.Range("J132").Value = 5344.4287
$ws.Range("K132").Value = 14243.2311
$ws.Range("L132").Value = 16033.2861
$ws.Range("M132").Value = -11713.2311
$ws.Range("N132").Value = -21093.2861

$ws.Range("H138").Value = 1466.2745
$ws.Range("I138").Value = 1081
$ws.Range("J138").Value = 2867.2727
$ws.Range("K138").Value = 3243
$ws.Range("L138").Value = 8601.8181
$ws.Range("M138").Value = 1897
$ws.Range("N138").Value = -18881.8181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6965.52
$ws.Range("I32").Value = 4603.827
$ws.Range("J32").Value = 17033.79
$ws.Range("K32").Value = 4603.827
$ws.Range("L32").Value = 17033.79
$ws.Range("M32").Value = -4316.827
$ws.Range("N32").Value = -17607.79

$ws.Range("H61").Value = 1561.0238
$ws.Range("I61").Value = 1388.68
$ws.Range("J61").Value = 1814.4706
$ws.Range("K61").Value = 1388.68
$ws.Range("L61").Value = 1814.4706
$ws.Range("M61").Value = -1176.68
$ws.Range("N61").Value = -2238.4706

$ws.Range("H122").Value = 1165.0555
$ws.Range("I122").Value = 879.2143
$ws.Range("J122").Value = 2165.5
$ws.Range("K122").Value = 2637.6429
$ws.Range("L122").Value = 6496.5
$ws.Range("M122").Value = -187.6428999999998
$ws.Range("N122").Value = -11396.5

$ws.Range("H132").Value = 1439965.2
$ws.Range("I132").Value = 2211.5112
$ws.Range("J132").Value = 2980415.8
$ws.Range("K132").Value = 6634.5336
$ws.Range("L132").Value = 8941247.399999999
$ws.Range("M132").Value = -4104.5336
$ws.Range("N132").Value = -8946307.399999999

$ws.Range("H136").Value = 1561.0238
$ws.Range("I136").Value = 1388.68
$ws.Range("J136").Value = 1814.4706
$ws.Range("K136").Value = 4166.04
$ws.Range("L136").Value = 5443.4118
$ws.Range("M136").Value = -1616.04
$ws.Range("N136").Value = -10543.4118

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1304.1428
$ws.Range("I107").Value = 1376
$ws.Range("J107").Value = 998.75
$ws.Range("K107").Value = 1376
$ws.Range("L107").Value = 998.75
$ws.Range("M107").Value = 544
$ws.Range("N107").Value = -4838.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2241.524
$ws.Range("I16").Value = 2205.0625
$ws.Range("J16").Value = 2358.2
$ws.Range("K16").Value = 2205.0625
$ws.Range("L16").Value = 2358.2
$ws.Range("M16").Value = -1918.0625
$ws.Range("N16").Value = -2932.2

$ws.Range("H31").Value = 4976932.5
$ws.Range("I31").Value = 1248.225
$ws.Range("J31").Value = 12348317
$ws.Range("K31").Value = 1248.225
$ws.Range("L31").Value = 12348317
$ws.Range("M31").Value = -953.2249999999999
$ws.Range("N31").Value = -12348907

$ws.Range("H34").Value = 4976932.5
$ws.Range("I34").Value = 1248.225
$ws.Range("J34").Value = 12348317
$ws.Range("K34").Value = 1248.225
$ws.Range("L34").Value = 12348317
$ws.Range("M34").Value = -1046.225
$ws.Range("N34").Value = -12348721

$ws.Range("H107").Value = 294.27274
$ws.Range("I107").Value = 270.77777
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 270.77777
$ws.Range("L107").Value = 400
$ws.Range("M107").Value = 1649.22223
$ws.Range("N107").Value = -4240

$ws.Range("H110").Value = 40000
$ws.Range("J110").Value = 40000
$ws.Range("L110").Value = 40000
$ws.Range("N110").Value = -48180

$ws.Range("H113").Value = 2241.524
$ws.Range("I113").Value = 2205.0625
$ws.Range("J113").Value = 2358.2
$ws.Range("K113").Value = 2205.0625
$ws.Range("L113").Value = 2358.2
$ws.Range("M113").Value = -35.0625
$ws.Range("N113").Value = -6698.2

$ws.Range("H132").Value = 2942.7878
$ws.Range("I132").Value = 2310.6155
$ws.Range("J132").Value = 3353.7
$ws.Range("K132").Value = 6931.8465
$ws.Range("L132").Value = 10061.1
$ws.Range("M132").Value = -4401.8465
$ws.Range("N132").Value = -15121.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1331.5555
$ws.Range("I113").Value = 496.25
$ws.Range("J113").Value = 1999.8
$ws.Range("K113").Value = 1488.75
$ws.Range("L113").Value = 5999.4
$ws.Range("M113").Value = 681.25
$ws.Range("N113").Value = -10339.4

$ws.Range("H114").Value = 2120.4167
$ws.Range("I114").Value = 189.5
$ws.Range("J114").Value = 2506.6
$ws.Range("K114").Value = 568.5
$ws.Range("L114").Value = 7519.799999999999
$ws.Range("M114").Value = 2685.5
$ws.Range("N114").Value = -14027.8

$ws.Range("H120").Value = 11381.111
$ws.Range("I120").Value = 9305
$ws.Range("J120").Value = 12419.167
$ws.Range("K120").Value = 27915
$ws.Range("L120").Value = 37257.501
$ws.Range("M120").Value = -23077
$ws.Range("N120").Value = -46933.501

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6330.864
$ws.Range("I70").Value = 5263.636
$ws.Range("J70").Value = 7398.091
$ws.Range("K70").Value = 5263.636
$ws.Range("L70").Value = 7398.091
$ws.Range("M70").Value = -4993.636
$ws.Range("N70").Value = -7938.091

$ws.Range("H73").Value = 6330.864
$ws.Range("I73").Value = 5263.636
$ws.Range("J73").Value = 7398.091
$ws.Range("K73").Value = 5263.636
$ws.Range("L73").Value = 7398.091
$ws.Range("M73").Value = -4327.636
$ws.Range("N73").Value = -9270.091

$ws.Range("H113").Value = 1066.1
$ws.Range("I113").Value = 1017.8889
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1017.8889
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 1152.1111
$ws.Range("N113").Value = -5840

$ws.Range("H122").Value = 2295.0588
$ws.Range("I122").Value = 2228.2964
$ws.Range("J122").Value = 2552.5715
$ws.Range("K122").Value = 6684.889200000001
$ws.Range("L122").Value = 7657.7145
$ws.Range("M122").Value = -4234.889200000001
$ws.Range("N122").Value = -12557.7145

$ws.Range("H126").Value = 12181.292
$ws.Range("I126").Value = 10681.1875
$ws.Range("J126").Value = 15181.5
$ws.Range("K126").Value = 32043.5625
$ws.Range("L126").Value = 45544.5
$ws.Range("M126").Value = -29573.5625
$ws.Range("N126").Value = -50484.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 11592.223
$ws.Range("J94").Value = 11592.223
$ws.Range("L94").Value = 11592.223
$ws.Range("N94").Value = -12944.223

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3189.8572
$ws.Range("I126").Value = 3600.25
$ws.Range("J126").Value = 727.5
$ws.Range("K126").Value = 10800.75
$ws.Range("L126").Value = 2182.5
$ws.Range("M126").Value = -8330.75
$ws.Range("N126").Value = -7122.5
